$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2021
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2021

$ws.Range("H2").Value = "ZZZ495"
$ws.Range("I2").Value = "ABC0987AX291"
$ws.Range("J2").Value = "MMAA09XFGS286"
$ws.Range("H3").Value = "ZZZ496"
$ws.Range("I3").Value = "ABC0987AX292"
$ws.Range("J3").Value = "MMAA09XFGS287"

$ws.Range("H4:J5").Clear()

$ws.Range("E4").Select()
